$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.982.84"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "'1.559.90"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'208.32"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'22.15"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "'0.249"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "'0.0599"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").Value = "'0.0856"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "'1.784.29"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "'1.564.33"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'3.76"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'61.94"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'26.982.33"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "'0.0₃0706"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "'216.78"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'7.37"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'4.12"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "'9.27"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").Value = "'153.14"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'6.61"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "'15.12"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D30").Value = "'0.0473"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "'3.18"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").Value = "'1.425.01"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("E36").Value = "  +7.84%  "
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").Value = "'0.0165"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("D39").Value = "'0.533"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "'5.85"
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").Value = "'64.59"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'1.75"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'1.695.85"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "'87.14"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").Value = "'0.0960"
$ws.Range("E51").Value = "  -0.05%  "
